# Rename the inline picture shapes embedded in the document's headers and
# footers. Word keeps these OOXML picture elements in sync via the
# InlineShape.Name property when a picture is renamed (e.g. after it is
# replaced/relinked in the Pearson/BTEC template), so we drive the change
# the same way an author renaming the assets in the UI would.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers -------------------------------------------------------------
# Primary footer ("default") -> previously named image1.png, becomes image2.png
$footerDefault = $sec.Footers.Item(1)
if ($footerDefault.Exists -and $footerDefault.Range.InlineShapes.Count -ge 1) {
    $footerDefault.Range.InlineShapes.Item(1).Name = "image2.png"
}

# First-page footer -> previously named image1.png, becomes image2.png
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $footerFirst.Range.InlineShapes.Item(1).Name = "image2.png"
}

# --- Headers ---------------------------------------------------------------
# First-page header (BTec logo) -> previously named image2.jpg, becomes image1.jpg
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $headerFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"
}
